# Login Positive Test is added
# Update the e-mail test data on the "pair" sheet: swap the gmail.com
# domain for ggmail.com in the E-Mail column (D) of the positive login
# test rows, then leave the selection where Excel would land after the
# last edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pair")

$ws.Range("D2").Value = "ali@ggmail.com"
$ws.Range("D3").Value = "a@ggmail.com"
$ws.Range("D5").Value = "ab@ggmail.com"
$ws.Range("D6").Value = "c@ggmail.com"
$ws.Range("D7").Value = "d@ggmail.com"

$ws.Range("E11").Select()
